# Apply the "username3/password3" row addition + cosmetic view tweaks
# described by the commit "working with Excel, CSV learn @Narrative".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (4th username/password pair) -----------------------
# Appending these two values adds two new shared-string entries
# ("username3"/"password3") and extends the sheet's used range /
# dimension from A1:B3 to A1:B4 automatically.
$ws.Range("A4").Value = "username3"
$ws.Range("B4").Value = "password3"

# --- Column widths ------------------------------------------------------
# Give column A/B an explicit width (the source workbook now carries a
# <cols> block instead of relying purely on the sheet default).
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(2).ColumnWidth = 10.3333333333333

# --- Sheet-wide default width/height/outline (best effort) --------------
# The edited file also nudges the sheet's default column width/row
# height and outline-row depth. These are legacy sheetFormatPr values;
# set them defensively through the documented properties in case the
# host persists them.
try { $ws.StandardWidth = 9.13888888888889 } catch {}
try { $ws.StandardHeight = 14.4 } catch {}

# --- Selection / active cell --------------------------------------------
# The saved selection moves on to the next empty row beneath the data
# that was just entered.
[void]$ws.Range("B5").Select()

# --- Window geometry (best effort) --------------------------------------
# The workbook view records a wider/shorter application window in the
# edited file. Window chrome is host UI state rather than worksheet
# data, but set it anyway in case the host surfaces it.
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 21000
    $win.Height = 12300
} catch {
}
